$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.1552973985672
$ws.Range("B1").Value = 2.322333335876465
$ws.Range("C1").Value = 4.408517360687256
$ws.Range("D1").Value = 3.524778127670288
$ws.Range("E1").Value = 1.236665725708008
